# Updates numeric stats (currentAveragePrice / LevePrice / LeveProfit columns)
# across the eight crafting-leve sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# to reflect the latest Market Board averages pulled by the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1657.7142
$ws.Range("J17").Value = 1657.7142
$ws.Range("L17").Value = 4973.142599999999
$ws.Range("N17").Value = -5309.142599999999
$ws.Range("H41").Value = 1775.1765
$ws.Range("I41").Value = 1607.9
$ws.Range("J41").Value = 2014.1428
$ws.Range("K41").Value = 1607.9
$ws.Range("L41").Value = 2014.1428
$ws.Range("M41").Value = -1167.9
$ws.Range("N41").Value = -2894.1428
$ws.Range("H69").Value = 17445.889
$ws.Range("I69").Value = 11285.571
$ws.Range("J69").Value = 39007
$ws.Range("K69").Value = 33856.713
$ws.Range("L69").Value = 117021
$ws.Range("M69").Value = -32982.713
$ws.Range("N69").Value = -118769
$ws.Range("H72").Value = 17445.889
$ws.Range("I72").Value = 11285.571
$ws.Range("J72").Value = 39007
$ws.Range("K72").Value = 101570.139
$ws.Range("L72").Value = 351063
$ws.Range("M72").Value = -97202.139
$ws.Range("N72").Value = -359799
$ws.Range("H106").Value = 8941.556
$ws.Range("I106").Value = 2268
$ws.Range("K106").Value = 2268
$ws.Range("M106").Value = -1637
$ws.Range("H111").Value = 1856.375
$ws.Range("I111").Value = 1332.2858
$ws.Range("K111").Value = 3996.8574
$ws.Range("M111").Value = -929.8574000000003
$ws.Range("H135").Value = 1162.5555
$ws.Range("I135").Value = 1205.6666
$ws.Range("K135").Value = 10850.9994
$ws.Range("M135").Value = -8315.999400000001
$ws.Range("H138").Value = 3189.6885
$ws.Range("I138").Value = 1198.75
$ws.Range("J138").Value = 3329.4036
$ws.Range("K138").Value = 3596.25
$ws.Range("L138").Value = 9988.210800000001
$ws.Range("M138").Value = 1543.75
$ws.Range("N138").Value = -20268.2108

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 1216.625
$ws.Range("I4").Value = 1233.75
$ws.Range("J4").Value = 1199.5
$ws.Range("K4").Value = 1233.75
$ws.Range("L4").Value = 1199.5
$ws.Range("M4").Value = -1117.75
$ws.Range("N4").Value = -1431.5
$ws.Range("H32").Value = 9806149
$ws.Range("I32").Value = 10870991
$ws.Range("K32").Value = 10870991
$ws.Range("M32").Value = -10870704
$ws.Range("H74").Value = 8072122.5
$ws.Range("I74").Value = 11365341
$ws.Range("J74").Value = 22032.223
$ws.Range("K74").Value = 11365341
$ws.Range("L74").Value = 22032.223
$ws.Range("M74").Value = -11364467
$ws.Range("N74").Value = -23780.223
$ws.Range("H77").Value = 8072122.5
$ws.Range("I77").Value = 11365341
$ws.Range("J77").Value = 22032.223
$ws.Range("K77").Value = 56826705
$ws.Range("L77").Value = 110161.115
$ws.Range("M77").Value = -56822337
$ws.Range("N77").Value = -118897.115

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3509.8823
$ws.Range("I20").Value = 3573
$ws.Range("J20").Value = 2500
$ws.Range("K20").Value = 3573
$ws.Range("L20").Value = 2500
$ws.Range("M20").Value = -3326
$ws.Range("N20").Value = -2994
$ws.Range("H22").Value = 390.8
$ws.Range("I22").Value = 113.5
$ws.Range("J22").Value = 1500
$ws.Range("K22").Value = 113.5
$ws.Range("L22").Value = 1500
$ws.Range("M22").Value = 59.5
$ws.Range("N22").Value = -1846
$ws.Range("H81").Value = 49721.5
$ws.Range("J81").Value = 49721.5
$ws.Range("L81").Value = 49721.5
$ws.Range("N81").Value = -51843.5
$ws.Range("H84").Value = 49721.5
$ws.Range("J84").Value = 49721.5
$ws.Range("L84").Value = 149164.5
$ws.Range("N84").Value = -159772.5
$ws.Range("H94").Value = 1043.04
$ws.Range("I94").Value = 1169
$ws.Range("J94").Value = 775.375
$ws.Range("K94").Value = 1169
$ws.Range("L94").Value = 775.375
$ws.Range("M94").Value = -718
$ws.Range("N94").Value = -1677.375
$ws.Range("H132").Value = 93888.89
$ws.Range("J132").Value = 83500
$ws.Range("L132").Value = 83500
$ws.Range("N132").Value = -93620
$ws.Range("H134").Value = 107301.5
$ws.Range("I134").Value = 1855.2858
$ws.Range("K134").Value = 5565.857400000001
$ws.Range("M134").Value = -3030.857400000001
$ws.Range("H140").Value = 80000
$ws.Range("J140").Value = 80000
$ws.Range("L140").Value = 80000
$ws.Range("N140").Value = -90360

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H31").Value = 408030.8
$ws.Range("I31").Value = 5042.6665
$ws.Range("K31").Value = 5042.6665
$ws.Range("M31").Value = -4747.6665
$ws.Range("H34").Value = 408030.8
$ws.Range("I34").Value = 5042.6665
$ws.Range("K34").Value = 5042.6665
$ws.Range("M34").Value = -4840.6665
$ws.Range("H106").Value = 38318.875
$ws.Range("J106").Value = 38318.875
$ws.Range("L106").Value = 38318.875
$ws.Range("N106").Value = -40842.875
$ws.Range("H132").Value = 2098.2778
$ws.Range("I132").Value = 1924.6
$ws.Range("K132").Value = 5773.799999999999
$ws.Range("M132").Value = -3243.799999999999
$ws.Range("H134").Value = 1003734.1
$ws.Range("I134").Value = 1113036.4
$ws.Range("K134").Value = 3339109.2
$ws.Range("M134").Value = -3336574.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1000
$ws.Range("I5").Value = 1000
$ws.Range("K5").Value = 3000
$ws.Range("M5").Value = -2888
$ws.Range("H23").Value = 230.875
$ws.Range("I23").Value = 150
$ws.Range("J23").Value = 279.4
$ws.Range("K23").Value = 450
$ws.Range("L23").Value = 838.1999999999999
$ws.Range("M23").Value = -215
$ws.Range("N23").Value = -1308.2
$ws.Range("H59").Value = 1625
$ws.Range("J59").Value = 2150
$ws.Range("L59").Value = 6450
$ws.Range("N59").Value = -7530
$ws.Range("H132").Value = 1805.7368
$ws.Range("I132").Value = 1964.3334
$ws.Range("K132").Value = 17679.0006
$ws.Range("M132").Value = -15149.0006
$ws.Range("H135").Value = 1000
$ws.Range("I135").Value = 1000
$ws.Range("K135").Value = 9000
$ws.Range("M135").Value = -6465

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4911.385
$ws.Range("I70").Value = 4820.6665
$ws.Range("K70").Value = 4820.6665
$ws.Range("M70").Value = -4550.6665
$ws.Range("H73").Value = 4911.385
$ws.Range("I73").Value = 4820.6665
$ws.Range("K73").Value = 4820.6665
$ws.Range("M73").Value = -3884.6665
$ws.Range("H104").Value = 47500
$ws.Range("J104").Value = 47500
$ws.Range("L104").Value = 47500
$ws.Range("N104").Value = -54488
$ws.Range("H122").Value = 1108.75
$ws.Range("I122").Value = 880.1667
$ws.Range("J122").Value = 1794.5
$ws.Range("K122").Value = 2640.5001
$ws.Range("L122").Value = 5383.5
$ws.Range("M122").Value = -190.5001000000002
$ws.Range("N122").Value = -10283.5
$ws.Range("H128").Value = 84663.336
$ws.Range("J128").Value = 84663.336
$ws.Range("L128").Value = 84663.336
$ws.Range("N128").Value = -94623.336
$ws.Range("H132").Value = 142862400
$ws.Range("I132").Value = 200006560
$ws.Range("K132").Value = 600019680
$ws.Range("M132").Value = -600017150

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("N38").ClearContents()
$ws.Range("H68").Value = 1773.75
$ws.Range("I68").Value = 1773.75
$ws.Range("K68").Value = 1773.75
$ws.Range("M68").Value = -1024.75
$ws.Range("H71").Value = 1773.75
$ws.Range("I71").Value = 1773.75
$ws.Range("K71").Value = 8868.75
$ws.Range("M71").Value = -5124.75
$ws.Range("H136").Value = 58645.26
$ws.Range("I136").Value = 9264.789000000001
$ws.Range("K136").Value = 27794.367
$ws.Range("M136").Value = -25244.367

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3241.1428
$ws.Range("I81").Value = 737.6
$ws.Range("J81").Value = 9500
$ws.Range("K81").Value = 1475.2
$ws.Range("L81").Value = 19000
$ws.Range("M81").Value = -414.2
$ws.Range("N81").Value = -21122
$ws.Range("H84").Value = 3241.1428
$ws.Range("I84").Value = 737.6
$ws.Range("J84").Value = 9500
$ws.Range("K84").Value = 7376
$ws.Range("L84").Value = 95000
$ws.Range("M84").Value = -2072
$ws.Range("N84").Value = -105608
$ws.Range("H100").Value = 1407.4
$ws.Range("I100").Value = 1484.25
$ws.Range("J100").Value = 1100
$ws.Range("K100").Value = 2968.5
$ws.Range("L100").Value = 2200
$ws.Range("M100").Value = -2427.5
$ws.Range("N100").Value = -3282

